# Auto-generated Excel COM-interop script applying numeric updates
# to the Maduin_Profits workbook, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1899.6666
$ws.Range("I32").Value = 1966.6666
$ws.Range("J32").Value = 1832.6666
$ws.Range("K32").Value = 1966.6666
$ws.Range("L32").Value = 1832.6666
$ws.Range("M32").Value = -1640.6666
$ws.Range("N32").Value = -2484.6666

$ws.Range("H48").Value = 1550
$ws.Range("I48").Value = 1500
$ws.Range("J48").Value = 1600
$ws.Range("K48").Value = 4500
$ws.Range("L48").Value = 4800
$ws.Range("M48").Value = -4208
$ws.Range("N48").Value = -5384

$ws.Range("H56").Value = 1550
$ws.Range("I56").Value = 1500
$ws.Range("J56").Value = 1600
$ws.Range("K56").Value = 4500
$ws.Range("L56").Value = 4800
$ws.Range("M56").Value = -3966
$ws.Range("N56").Value = -5868

$ws.Range("H64").Value = 14443.444
$ws.Range("I64").Value = 8888.111000000001
$ws.Range("J64").Value = 19998.777
$ws.Range("K64").Value = 8888.111000000001
$ws.Range("L64").Value = 19998.777
$ws.Range("M64").Value = -8640.111000000001
$ws.Range("N64").Value = -20494.777

$ws.Range("H67").Value = 14443.444
$ws.Range("I67").Value = 8888.111000000001
$ws.Range("J67").Value = 19998.777
$ws.Range("K67").Value = 8888.111000000001
$ws.Range("L67").Value = 19998.777
$ws.Range("M67").Value = -8030.111000000001
$ws.Range("N67").Value = -21714.777

$ws.Range("H70").Value = 1150
$ws.Range("I70").Value = 1250
$ws.Range("J70").Value = 1050
$ws.Range("K70").Value = 3750
$ws.Range("L70").Value = 3150
$ws.Range("M70").Value = -3480
$ws.Range("N70").Value = -3690

$ws.Range("H73").Value = 1150
$ws.Range("I73").Value = 1250
$ws.Range("J73").Value = 1050
$ws.Range("K73").Value = 3750
$ws.Range("L73").Value = 3150
$ws.Range("M73").Value = -2814
$ws.Range("N73").Value = -5022

$ws.Range("H81").Value = 99999
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 99999
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 99999
$ws.Range("N81").Value = -101995
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 99999
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 99999
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 299997
$ws.Range("N84").Value = -309981
$ws.Range("M84").ClearContents()

$ws.Range("H100").Value = 1981.1786
$ws.Range("I100").Value = 1961
$ws.Range("J100").Value = 2102.25
$ws.Range("K100").Value = 1961
$ws.Range("L100").Value = 2102.25
$ws.Range("M100").Value = -1420
$ws.Range("N100").Value = -3184.25

$ws.Range("H113").Value = 1275.8182
$ws.Range("I113").Value = 1322.5
$ws.Range("J113").Value = 809
$ws.Range("K113").Value = 1322.5
$ws.Range("L113").Value = 809
$ws.Range("M113").Value = 1931.5
$ws.Range("N113").Value = -7317

$ws.Range("H137").Value = 2499.75
$ws.Range("I137").Value = 2499
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 7497
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -4947
$ws.Range("N137").Value = -12600

$ws.Range("H140").Value = 85999.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 85999.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 85999.8
$ws.Range("N140").Value = -96359.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2610.75
$ws.Range("I86").Value = 2649.3333
$ws.Range("J86").Value = 2495
$ws.Range("K86").Value = 2649.3333
$ws.Range("L86").Value = 2495
$ws.Range("M86").Value = -1526.3333
$ws.Range("N86").Value = -4741

$ws.Range("H89").Value = 2610.75
$ws.Range("I89").Value = 2649.3333
$ws.Range("J89").Value = 2495
$ws.Range("K89").Value = 13246.6665
$ws.Range("L89").Value = 12475
$ws.Range("M89").Value = -7630.666499999999
$ws.Range("N89").Value = -23707

$ws.Range("H99").Value = 1894.2307
$ws.Range("I99").Value = 1428.125
$ws.Range("J99").Value = 2640
$ws.Range("K99").Value = 1428.125
$ws.Range("L99").Value = 2640
$ws.Range("M99").Value = 69.875
$ws.Range("N99").Value = -5636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1031.2858
$ws.Range("I16").Value = 807.75
$ws.Range("J16").Value = 1329.3334
$ws.Range("K16").Value = 807.75
$ws.Range("L16").Value = 1329.3334
$ws.Range("M16").Value = -520.75
$ws.Range("N16").Value = -1903.3334

$ws.Range("H18").Value = 46242.125
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 46242.125
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 46242.125
$ws.Range("N18").Value = -46702.125

$ws.Range("H31").Value = 3201.8333
$ws.Range("I31").Value = 2973
$ws.Range("J31").Value = 3888.3333
$ws.Range("K31").Value = 2973
$ws.Range("L31").Value = 3888.3333
$ws.Range("M31").Value = -2678
$ws.Range("N31").Value = -4478.3333

$ws.Range("H34").Value = 3201.8333
$ws.Range("I34").Value = 2973
$ws.Range("J34").Value = 3888.3333
$ws.Range("K34").Value = 2973
$ws.Range("L34").Value = 3888.3333
$ws.Range("M34").Value = -2771
$ws.Range("N34").Value = -4292.3333

$ws.Range("H43").Value = 39666.332
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 39666.332
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 39666.332
$ws.Range("N43").Value = -40034.332

$ws.Range("H69").Value = 7238.6924
$ws.Range("I69").Value = 5092
$ws.Range("J69").Value = 32999
$ws.Range("K69").Value = 5092
$ws.Range("L69").Value = 32999
$ws.Range("M69").Value = -4343
$ws.Range("N69").Value = -34497

$ws.Range("H72").Value = 7238.6924
$ws.Range("I72").Value = 5092
$ws.Range("J72").Value = 32999
$ws.Range("K72").Value = 15276
$ws.Range("L72").Value = 98997
$ws.Range("M72").Value = -11532
$ws.Range("N72").Value = -106485

$ws.Range("H99").Value = 6498.0527
$ws.Range("I99").Value = 5838.5713
$ws.Range("J99").Value = 8344.6
$ws.Range("K99").Value = 5838.5713
$ws.Range("L99").Value = 8344.6
$ws.Range("M99").Value = -4340.5713
$ws.Range("N99").Value = -11340.6

$ws.Range("H101").Value = 39666.332
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 39666.332
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 39666.332
$ws.Range("N101").Value = -46156.332

$ws.Range("H102").Value = 11000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 11000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 11000
$ws.Range("N102").Value = -15868

$ws.Range("H113").Value = 1031.2858
$ws.Range("I113").Value = 807.75
$ws.Range("J113").Value = 1329.3334
$ws.Range("K113").Value = 807.75
$ws.Range("L113").Value = 1329.3334
$ws.Range("M113").Value = 1362.25
$ws.Range("N113").Value = -5669.3334

$ws.Range("H126").Value = 6498.0527
$ws.Range("I126").Value = 5838.5713
$ws.Range("J126").Value = 8344.6
$ws.Range("K126").Value = 17515.7139
$ws.Range("L126").Value = 25033.8
$ws.Range("M126").Value = -15045.7139
$ws.Range("N126").Value = -29973.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 17796
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 17796
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 17796
$ws.Range("N15").Value = -18372

$ws.Range("H81").Value = 17796
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 17796
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 17796
$ws.Range("N81").Value = -19792

$ws.Range("H84").Value = 17796
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 17796
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 53388
$ws.Range("N84").Value = -63372

$ws.Range("H126").Value = 250007460
$ws.Range("I126").Value = 1000000000
$ws.Range("J126").Value = 9932.666999999999
$ws.Range("K126").Value = 3000000000
$ws.Range("L126").Value = 29798.001
$ws.Range("M126").Value = -2999997530
$ws.Range("N126").Value = -34738.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 32000
$ws.Range("I25").Value = 35000
$ws.Range("J25").Value = 26000
$ws.Range("K25").Value = 35000
$ws.Range("L25").Value = 26000
$ws.Range("M25").Value = -34770
$ws.Range("N25").Value = -26460

$ws.Range("H93").Value = 3500
$ws.Range("I93").Value = 3500
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3500
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2252

$ws.Range("H100").Value = 1499
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1499
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 1499
$ws.Range("N100").Value = -2581
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 25000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 25000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26622

$ws.Range("H71").Value = 25000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 25000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -83112
